# Pushed back sprints on Product Backlog
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Sprint #") values each pushed back by one sprint (mostly +1),
# row 4 became a merged-sprint label "2 and 3" (formatted like a date,
# d-mmm, matching the source workbook's cell format for that cell), and
# row 27 was pushed back by two sprints.
$ws.Range("C3").Value  = 10
$ws.Range("C5").Value  = 4
$ws.Range("C6").Value  = 5
$ws.Range("C7").Value  = 6
$ws.Range("C8").Value  = 5
$ws.Range("C9").Value  = 6
$ws.Range("C10").Value = 9
$ws.Range("C11").Value = 9
$ws.Range("C12").Value = 6
$ws.Range("C13").Value = 10
$ws.Range("C14").Value = 4
$ws.Range("C15").Value = 7
$ws.Range("C16").Value = 5
$ws.Range("C17").Value = 10
$ws.Range("C18").Value = 9
$ws.Range("C20").Value = 7
$ws.Range("C21").Value = 8
$ws.Range("C22").Value = 7
$ws.Range("C23").Value = 8
$ws.Range("C24").Value = 9
$ws.Range("C25").Value = 7
$ws.Range("C26").Value = 9
$ws.Range("C27").Value = 10
$ws.Range("C28").Value = 9
$ws.Range("C29").Value = 10

# Row 4 spans two sprints, so it gets a text label instead of a number.
$ws.Range("C4").NumberFormat = "d-mmm"
$ws.Range("C4").Value = "2 and 3"

# Leave the final selection on C29, matching the saved workbook state.
$ws.Range("C29").Select()
